$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 45, pushing existing rows 45-55 down to 47-57.
$ws.Rows.Item(45).Insert()
$ws.Rows.Item(45).Insert()

# New row 45: weekly "Especial" quality entry
$ws.Cells.Item(45, 1).Value = 11
$ws.Cells.Item(45, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(45, 3).Value = "Bíobío"
$ws.Cells.Item(45, 4).Value = 45218
$ws.Cells.Item(45, 5).Value = 8
$ws.Cells.Item(45, 6).Value = "Fruta"
$ws.Cells.Item(45, 7).Value = 100107
$ws.Cells.Item(45, 8).Value = "Otros"
$ws.Cells.Item(45, 9).Value = 100107002
$ws.Cells.Item(45, 10).Value = "Chirimoya"
$ws.Cells.Item(45, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(45, 12).Value = "Especial"
$ws.Cells.Item(45, 13).Value = 100
$ws.Cells.Item(45, 14).Value = 20000
$ws.Cells.Item(45, 15).Value = 20000
$ws.Cells.Item(45, 16).Value = 20000
$ws.Cells.Item(45, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(45, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(45, 19).Value = 2000
$ws.Cells.Item(45, 20).Value = 10

# New row 46: weekly "Primera" quality entry
$ws.Cells.Item(46, 1).Value = 11
$ws.Cells.Item(46, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(46, 3).Value = "Bíobío"
$ws.Cells.Item(46, 4).Value = 45218
$ws.Cells.Item(46, 5).Value = 8
$ws.Cells.Item(46, 6).Value = "Fruta"
$ws.Cells.Item(46, 7).Value = 100107
$ws.Cells.Item(46, 8).Value = "Otros"
$ws.Cells.Item(46, 9).Value = 100107002
$ws.Cells.Item(46, 10).Value = "Chirimoya"
$ws.Cells.Item(46, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(46, 12).Value = "Primera"
$ws.Cells.Item(46, 13).Value = 80
$ws.Cells.Item(46, 14).Value = 17000
$ws.Cells.Item(46, 15).Value = 17000
$ws.Cells.Item(46, 16).Value = 17000
$ws.Cells.Item(46, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(46, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(46, 19).Value = 1700
$ws.Cells.Item(46, 20).Value = 10
